$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "64.383.94"
$ws.Range("E2").Value2 = "  +0.22%  "
$ws.Range("D3").Value2 = "3.505.89"
$ws.Range("E3").Value2 = "  +0.05%  "
$ws.Range("E4").Value2 = "  -0.07%  "
$ws.Range("D5").Value2 = "'590.09"
$ws.Range("E5").Value2 = "  +1.07%  "
$ws.Range("D6").Value2 = "'134.46"
$ws.Range("E6").Value2 = "  -0.13%  "
$ws.Range("E8").Value2 = "  -0.10%  "
$ws.Range("E9").Value2 = "  +5.25%  "
$ws.Range("E10").Value2 = "  -0.02%  "
$ws.Range("E11").Value2 = "  +2.19%  "
$ws.Range("D12").Value2 = "4.101.94"
$ws.Range("E12").Value2 = "  -0.04%  "
$ws.Range("E13").Value2 = "  +1.61%  "
$ws.Range("E14").Value2 = "  +1.00%  "
$ws.Range("D15").Value2 = "3.503.96"
$ws.Range("E15").Value2 = "  +0.03%  "
$ws.Range("D16").Value2 = "'25.78"
$ws.Range("E16").Value2 = "  -5.54%  "
$ws.Range("D17").Value2 = "64.366.14"
$ws.Range("E17").Value2 = "  +0.14%  "
$ws.Range("D18").Value2 = "'9.91"
$ws.Range("E18").Value2 = "  +1.16%  "
$ws.Range("D19").Value2 = "'5.76"
$ws.Range("E19").Value2 = "  +2.93%  "
$ws.Range("E20").Value2 = "  -2.26%  "
$ws.Range("D21").Value2 = "'393.36"
$ws.Range("E21").Value2 = "  +2.48%  "
$ws.Range("E22").Value2 = "  +1.18%  "
$ws.Range("D23").Value2 = "3.645.88"
$ws.Range("E23").Value2 = "  +0.01%  "
$ws.Range("D24").Value2 = "'74.65"
$ws.Range("E24").Value2 = "  +1.05%  "
$ws.Range("E25").Value2 = "  +0.15%  "
$ws.Range("E26").Value2 = "  +0.11%  "
$ws.Range("E27").Value2 = "  +2.26%  "
$ws.Range("E28").Value2 = "  +0.09%  "
$ws.Range("D29").Value2 = "'7.39"
$ws.Range("E29").Value2 = "  -3.00%  "
$ws.Range("E30").Value2 = "  +2.07%  "
$ws.Range("E31").Value2 = "  -0.81%  "
$ws.Range("E32").Value2 = "  -6.94%  "
$ws.Range("E33").Value2 = "  +8.25%  "
$ws.Range("D34").Value2 = "3.530.86"
$ws.Range("E34").Value2 = "  +0.33%  "
$ws.Range("E35").Value2 = "  +0.02%  "
$ws.Range("E36").Value2 = "  -0.87%  "
$ws.Range("D37").Value2 = "'5.34"
$ws.Range("E37").Value2 = "  +0.71%  "
$ws.Range("E38").Value2 = "  +0.99%  "
$ws.Range("D39").Value2 = "'1.55"
$ws.Range("E39").Value2 = "  -0.32%  "
$ws.Range("D40").Value2 = "'167.52"
$ws.Range("E40").Value2 = "  +2.11%  "
$ws.Range("D41").Value2 = "'0.0788"
$ws.Range("E41").Value2 = "  +0.25%  "
$ws.Range("E42").Value2 = "  +0.11%  "
$ws.Range("E43").Value2 = "  -0.10%  "
$ws.Range("D44").Value2 = "'4.45"
$ws.Range("E44").Value2 = "  +1.21%  "
$ws.Range("D45").Value2 = "'24.92"
$ws.Range("E45").Value2 = "  -4.63%  "
$ws.Range("E46").Value2 = "  +2.31%  "
$ws.Range("D47").Value2 = "'1.18"
$ws.Range("E47").Value2 = "  -3.54%  "
$ws.Range("E48").Value2 = "  +0.50%  "
$ws.Range("D49").Value2 = "2.386.08"
$ws.Range("E49").Value2 = "  -3.73%  "
$ws.Range("E50").Value2 = "  -1.91%  "
$ws.Range("B51").Value2 = "TheGraph"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").Value2 = "'0.219"
$ws.Range("E51").Value2 = "  +1.09%  "
